{"js": "// Update the \"Skills\" section: each bullet's bold category label and its\n// following description are replaced (and several categories are renamed\n// / reordered / expanded) per the resume refresh.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Locate the \"Skills\" Heading2 paragraph; the 8 rows that follow it are the\n// skill category lines this edit touches.\nlet skillsIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Skills\") {\n    skillsIndex = i;\n    break;\n  }\n}\nif (skillsIndex === -1) {\n  throw new Error('Could not find the \"Skills\" heading paragraph.');\n}\n\n// [old bold label, new bold label, new description (leading space, no colon)]\nconst updates = [\n  [\n    \"Backend Development\",\n    \"Core Programming & Development\",\n    \" C#, .NET, .NET Core, .NET Framework, .NET Standard, .NET 6, .NET 8, LINQ, TypeScript, JavaScript, REST APIs, RESTful APIs, JSON, XML, Webhooks, PowerShell, Bash, VB.NET, Code Generation, Windows Services, xUnit, NUnit, NCover, Automated Testing\",\n  ],\n  [\n    \"Frontend Development\",\n    \"Frontend Technologies\",\n    \" React, Angular, AngularJS, TypeScript, JavaScript, HTML, CSS, Tailwind CSS, Bootstrap, LESS, jQuery, AJAX, Schema.org, SEO Optimization, Web Development, Google Analytics, Google Webmaster Tools, Google AdSense\",\n  ],\n  [\n    \"Databases & Search\",\n    \"Backend & Server Technologies\",\n    \" ASP.NET MVC, ASP.NET Core, ASP.NET Razor, ASP.NET Web API, ASP.NET, Windows Services, API Integration, Authentication/Authorization, OAuth, IIS, Linux, Ubuntu, HAProxy, MSMQ, Event-Driven Architecture, .NET Web Services, .NET Remoting, BizTalk, Commerce Server, Classic ASP\",\n  ],\n  [\n    \"Cloud & DevOps\",\n    \"Data Management\",\n    \" SQL Server, T-SQL, PostgreSQL, MySQL, DynamoDB, Redis, Elasticsearch, OpenSearch, Couchbase, RavenDB, MongoDB, Lucene, Solr, Neo4j, Oracle, InfluxDB, AWS ElastiCache, AWS Redshift, AWS CloudSearch, Data Migration, Data Validation\",\n  ],\n  [\n    \"Architecture & Design\",\n    \"Cloud & Infrastructure\",\n    \" AWS, AWS Lambda, AWS S3, AWS CloudFront, AWS Elastic Beanstalk, AWS DynamoDB, AWS CloudWatch, AWS SES, AWS Route 53, AWS WAF, AWS EC2, AWS SNS, AWS EMR, AWS DataSync, Azure, Azure DevOps, CI/CD Pipelines, Git, GitHub, Jenkins, Bitbucket, Ansible, DevOps, Build Automation, CruiseControl.NET, MSBuild, NAnt, Web Deploy, Infrastructure Management, Nagios, NagiosXI, Kinesis\",\n  ],\n  [\n    \"AI & Machine Learning\",\n    \"Architecture & Performance\",\n    \" Software Architecture, Enterprise Architecture, Enterprise Systems Design, Event-Driven Architecture, Microservices, Scalable Systems, High-Volume Processing, Grid Computing, Performance Tuning, Software Development Life Cycle (SDLC), Code Generation, Process Improvement, Technical Documentation, Algorithm Design\",\n  ],\n  [\n    \"Blockchain & Web3\",\n    \"Specialized Technologies & Domains\",\n    \" OpenAI GPT, Machine Learning Integration, Confidence Scoring, HIVE Engine Blockchain, Web3, NFTs, Cryptocurrency, Decentralized Applications (dApps), Smart Contracts, Blockchain Technology, Game Development, Game Design, Game Monetization, Game Analytics, SharePoint, Documentum, DotNetNuke, PHP, Costpoint, ERP Systems, VB6, Visual Basic, ActiveX, COM, COM Interop\",\n  ],\n  [\n    \"Leadership\",\n    \"Leadership & Project Management\",\n    \" Team Leadership, Technical Leadership, Executive Leadership, Mentoring, Cross-team Coordination, Strategic Planning, Succession Planning, Requirements Gathering, Product Management, Agile, Scrum, Scrum Master, Release Management, Trello, Jira, Confluence, FishEye, Microsoft Teams\",\n  ],\n];\n\n// Each skill row is \"<bold label>: <description>\" stored as two runs. Split\n// each paragraph on its single colon to get a range for the bold label and a\n// range for the description, then replace each in place so the bold\n// formatting on the label run is preserved.\nconst rowParagraphs = [];\nfor (let r = 0; r < updates.length; r++) {\n  rowParagraphs.push(paragraphs.items[skillsIndex + 1 + r]);\n}\n\nconst rowParts = rowParagraphs.map((p) => p.split([\":\"], true, false));\nfor (const parts of rowParts) {\n  parts.load(\"items/text\");\n}\nawait context.sync();\n\nfor (let r = 0; r < updates.length; r++) {\n  const [oldLabel, newLabel, newDescription] = updates[r];\n  const parts = rowParts[r].items;\n  if (parts.length < 2) {\n    throw new Error(`Unexpected skill row format at row ${r}.`);\n  }\n  if (parts[0].text !== oldLabel) {\n    throw new Error(\n      `Skill row ${r} label mismatch: expected \"${oldLabel}\", found \"${parts[0].text}\".`\n    );\n  }\n  parts[0].insertText(newLabel, \"Replace\");\n  parts[1].insertText(newDescription, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Skills\" section: each bullet's bold category label and its\n# following description are replaced (and several categories are renamed /\n# reordered / expanded) per the resume refresh.\n$d = $word.ActiveDocument\n\n# Locate the \"Skills\" Heading2 paragraph; the 8 rows that follow it are the\n# skill category lines this edit touches.\n$paraCount = $d.Paragraphs.Count\n$skillsIndex = -1\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.TrimEnd(\"`r\") -eq \"Skills\") {\n        $skillsIndex = $i\n        break\n    }\n}\nif ($skillsIndex -eq -1) {\n    throw 'Could not find the \"Skills\" heading paragraph.'\n}\n\n# Each entry: old bold label, new bold label, new description (leading\n# space, no colon).\n$updates = @(\n    @(\"Backend Development\", \"Core Programming & Development\", \" C#, .NET, .NET Core, .NET Framework, .NET Standard, .NET 6, .NET 8, LINQ, TypeScript, JavaScript, REST APIs, RESTful APIs, JSON, XML, Webhooks, PowerShell, Bash, VB.NET, Code Generation, Windows Services, xUnit, NUnit, NCover, Automated Testing\"),\n    @(\"Frontend Development\", \"Frontend Technologies\", \" React, Angular, AngularJS, TypeScript, JavaScript, HTML, CSS, Tailwind CSS, Bootstrap, LESS, jQuery, AJAX, Schema.org, SEO Optimization, Web Development, Google Analytics, Google Webmaster Tools, Google AdSense\"),\n    @(\"Databases & Search\", \"Backend & Server Technologies\", \" ASP.NET MVC, ASP.NET Core, ASP.NET Razor, ASP.NET Web API, ASP.NET, Windows Services, API Integration, Authentication/Authorization, OAuth, IIS, Linux, Ubuntu, HAProxy, MSMQ, Event-Driven Architecture, .NET Web Services, .NET Remoting, BizTalk, Commerce Server, Classic ASP\"),\n    @(\"Cloud & DevOps\", \"Data Management\", \" SQL Server, T-SQL, PostgreSQL, MySQL, DynamoDB, Redis, Elasticsearch, OpenSearch, Couchbase, RavenDB, MongoDB, Lucene, Solr, Neo4j, Oracle, InfluxDB, AWS ElastiCache, AWS Redshift, AWS CloudSearch, Data Migration, Data Validation\"),\n    @(\"Architecture & Design\", \"Cloud & Infrastructure\", \" AWS, AWS Lambda, AWS S3, AWS CloudFront, AWS Elastic Beanstalk, AWS DynamoDB, AWS CloudWatch, AWS SES, AWS Route 53, AWS WAF, AWS EC2, AWS SNS, AWS EMR, AWS DataSync, Azure, Azure DevOps, CI/CD Pipelines, Git, GitHub, Jenkins, Bitbucket, Ansible, DevOps, Build Automation, CruiseControl.NET, MSBuild, NAnt, Web Deploy, Infrastructure Management, Nagios, NagiosXI, Kinesis\"),\n    @(\"AI & Machine Learning\", \"Architecture & Performance\", \" Software Architecture, Enterprise Architecture, Enterprise Systems Design, Event-Driven Architecture, Microservices, Scalable Systems, High-Volume Processing, Grid Computing, Performance Tuning, Software Development Life Cycle (SDLC), Code Generation, Process Improvement, Technical Documentation, Algorithm Design\"),\n    @(\"Blockchain & Web3\", \"Specialized Technologies & Domains\", \" OpenAI GPT, Machine Learning Integration, Confidence Scoring, HIVE Engine Blockchain, Web3, NFTs, Cryptocurrency, Decentralized Applications (dApps), Smart Contracts, Blockchain Technology, Game Development, Game Design, Game Monetization, Game Analytics, SharePoint, Documentum, DotNetNuke, PHP, Costpoint, ERP Systems, VB6, Visual Basic, ActiveX, COM, COM Interop\"),\n    @(\"Leadership\", \"Leadership & Project Management\", \" Team Leadership, Technical Leadership, Executive Leadership, Mentoring, Cross-team Coordination, Strategic Planning, Succession Planning, Requirements Gathering, Product Management, Agile, Scrum, Scrum Master, Release Management, Trello, Jira, Confluence, FishEye, Microsoft Teams\")\n)\n\n# Each skill row is \"<bold label>: <description>\" stored as two runs. Split\n# each paragraph's range on its single colon to get a sub-range for the bold\n# label and a sub-range for the description, then replace each in place\n# (working from the last row back to the first) so the bold formatting on\n# the label run is preserved and earlier offsets are unaffected by later\n# length changes.\nfor ($row = $updates.Count - 1; $row -ge 0; $row--) {\n    $oldLabel = $updates[$row][0]\n    $newLabel = $updates[$row][1]\n    $newDescription = $updates[$row][2]\n\n    $paraIndex = $skillsIndex + 1 + $row\n    $p = $d.Paragraphs.Item($paraIndex)\n    $r = $p.Range\n    $fullText = $r.Text.TrimEnd(\"`r\")\n    $colonIdx = $fullText.IndexOf(\":\")\n    if ($colonIdx -lt 0) {\n        throw \"Unexpected skill row format at paragraph $paraIndex.\"\n    }\n\n    $start = $r.Start\n    $labelRange = $d.Range($start, $start + $colonIdx)\n    if ($labelRange.Text -ne $oldLabel) {\n        throw \"Skill row at paragraph $paraIndex label mismatch: expected '$oldLabel', found '$($labelRange.Text)'.\"\n    }\n\n    $descRange = $d.Range($start + $colonIdx + 1, $r.End - 1)\n\n    $descRange.Text = $newDescription\n    $labelRange.Text = $newLabel\n}\n"}
